$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.162114024162292
$ws.Range("B1").Value = 2.582020282745361
$ws.Range("C1").Value = 4.202258110046387
$ws.Range("D1").Value = 3.443256139755249
$ws.Range("E1").Value = 1.21016788482666
